$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.595116972923279
$ws.Range("B1").Value = 3.197046995162964
$ws.Range("C1").Value = 2.902799844741821
$ws.Range("D1").Value = 1.777951002120972
$ws.Range("E1").Value = 0.9104089140892029
